$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "28.247.90"
$c.ClearFormats()
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "  -0.54%  "
$c.ClearFormats()

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.804.87"
$c.ClearFormats()
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "  -0.68%  "
$c.ClearFormats()

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = "  +0.06%  "
$c.ClearFormats()

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "314.69"
$c.ClearFormats()
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "  -0.11%  "
$c.ClearFormats()

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.ClearFormats()
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "  +0.10%  "
$c.ClearFormats()

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5269"
$c.ClearFormats()
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "  +2.71%  "
$c.ClearFormats()

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3825"
$c.ClearFormats()
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "  -2.75%  "
$c.ClearFormats()

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.08019"
$c.ClearFormats()
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "  +0.82%  "
$c.ClearFormats()

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "41.44"
$c.ClearFormats()
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "  -0.55%  "
$c.ClearFormats()

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "1.101"
$c.ClearFormats()
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "  -0.52%  "
$c.ClearFormats()

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "6.337"
$c.ClearFormats()
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "  +1.44%  "
$c.ClearFormats()

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "  +0.14%  "
$c.ClearFormats()

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "20.63"
$c.ClearFormats()
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "  -1.65%  "
$c.ClearFormats()

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.341"
$c.ClearFormats()
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "  -1.87%  "
$c.ClearFormats()

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.804.51"
$c.ClearFormats()
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "  -0.80%  "
$c.ClearFormats()

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "92.45"
$c.ClearFormats()
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "  -0.07%  "
$c.ClearFormats()

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.00001097"
$c.ClearFormats()
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "  -2.75%  "
$c.ClearFormats()

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06614"
$c.ClearFormats()
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "  -0.13%  "
$c.ClearFormats()

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "  +0.12%  "
$c.ClearFormats()

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "17.40"
$c.ClearFormats()
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "  -1.46%  "
$c.ClearFormats()

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.970"
$c.ClearFormats()
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "  -1.89%  "
$c.ClearFormats()

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "28.307.08"
$c.ClearFormats()
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "  -0.46%  "
$c.ClearFormats()

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "11.16"
$c.ClearFormats()
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "  -0.88%  "
$c.ClearFormats()

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.253"
$c.ClearFormats()
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "  -0.57%  "
$c.ClearFormats()

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "160.69"
$c.ClearFormats()
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "  +3.42%  "
$c.ClearFormats()

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = "  -2.89%  "
$c.ClearFormats()

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.010.37"
$c.ClearFormats()
$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = "  -0.86%  "
$c.ClearFormats()

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.362"
$c.ClearFormats()
$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = "  -1.61%  "
$c.ClearFormats()

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "123.32"
$c.ClearFormats()
$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = "  -1.72%  "
$c.ClearFormats()

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.1086"
$c.ClearFormats()
$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = "  -1.30%  "
$c.ClearFormats()

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.059"
$c.ClearFormats()
$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = "  -3.84%  "
$c.ClearFormats()

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.686"
$c.ClearFormats()
$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = "  +0.95%  "
$c.ClearFormats()

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.561"
$c.ClearFormats()
$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = "  -2.27%  "
$c.ClearFormats()

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.07292"
$c.ClearFormats()
$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = "  +3.87%  "
$c.ClearFormats()

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "12.48"
$c.ClearFormats()
$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = "  +10.56%  "
$c.ClearFormats()

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.02324"
$c.ClearFormats()
$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = "  +0.00%  "
$c.ClearFormats()

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2162"
$c.ClearFormats()
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "  -2.55%  "
$c.ClearFormats()

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "  -0.36%  "
$c.ClearFormats()

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "8.697"
$c.ClearFormats()
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "  -1.52%  "
$c.ClearFormats()

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.6227"
$c.ClearFormats()
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "  -0.32%  "
$c.ClearFormats()

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.169"
$c.ClearFormats()
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "  -0.48%  "
$c.ClearFormats()

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.372"
$c.ClearFormats()
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "  -2.05%  "
$c.ClearFormats()

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.6051"
$c.ClearFormats()
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "  +2.60%  "
$c.ClearFormats()

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "13.22"
$c.ClearFormats()
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "  -1.85%  "
$c.ClearFormats()

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.773"
$c.ClearFormats()
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "  +0.85%  "
$c.ClearFormats()

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "127.12"
$c.ClearFormats()
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "  +1.80%  "
$c.ClearFormats()

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.218"
$c.ClearFormats()
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "  +2.61%  "
$c.ClearFormats()

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.931"
$c.ClearFormats()
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "  -2.03%  "
$c.ClearFormats()

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06820"
$c.ClearFormats()
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "  -0.93%  "
$c.ClearFormats()

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "73.16"
$c.ClearFormats()
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "  -1.81%  "
$c.ClearFormats()
